$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (66) down into the new row (67)
$ws.Range("A66:F66").Copy()
$ws.Range("A67:F67").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new day's data (2020-05-18, serial date 43969)
$ws.Range("A67").Value2 = 43969
$ws.Range("B67").Value2 = 532
$ws.Range("C67").Value2 = 152
$ws.Range("D67").Value2 = 275
$ws.Range("E67").Value2 = 12
$ws.Range("F67").Value2 = 17

# Grow the table (ListObject) so the new row becomes part of it
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F67"))

# Match the post-edit selection state from the workbook
$ws.Range("E68").Select()
